# Weekly data refresh: a new week of observations (date 44585) was added
# to the "Apio" (Celery) price sheet. This pushes all existing data rows
# down by two rows (one new "Primera" row + one new "Segunda" row).
#
# The workbook has a single sheet with header row 1 and data rows 2..527.
# We insert two new blank rows right before the old row 494 (the start of
# the most-recent weekly pair) and populate them with the new week's data.
# Excel automatically shifts every subsequent row down by two and extends
# the used range / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 494, shifting existing rows 494:527 down
# to become rows 496:529.
$ws.Rows("494:495").Insert()

# New row 494 - "Primera" quality observation for the new week (44585).
$ws.Cells.Item(494, 1).Value = 6
$ws.Cells.Item(494, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(494, 3).Value = 'Metropolitana'
$ws.Cells.Item(494, 4).Value = 44585
$ws.Cells.Item(494, 5).Value = 13
$ws.Cells.Item(494, 6).Value = 100112017
$ws.Cells.Item(494, 7).Value = 'Apio'
$ws.Cells.Item(494, 8).Value = 'Americana (o)'
$ws.Cells.Item(494, 9).Value = 'Primera'
$ws.Cells.Item(494, 10).Value = 2300
$ws.Cells.Item(494, 11).Value = 6000
$ws.Cells.Item(494, 12).Value = 7000
$ws.Cells.Item(494, 13).Value = 6391
$ws.Cells.Item(494, 14).Value = '$/docena de matas'
$ws.Cells.Item(494, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(494, 16).Value = 1065
$ws.Cells.Item(494, 17).Value = 6
$ws.Cells.Item(494, 18).Value = 'Hortaliza'

# New row 495 - "Segunda" quality observation for the new week (44585).
$ws.Cells.Item(495, 1).Value = 6
$ws.Cells.Item(495, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(495, 3).Value = 'Metropolitana'
$ws.Cells.Item(495, 4).Value = 44585
$ws.Cells.Item(495, 5).Value = 13
$ws.Cells.Item(495, 6).Value = 100112017
$ws.Cells.Item(495, 7).Value = 'Apio'
$ws.Cells.Item(495, 8).Value = 'Americana (o)'
$ws.Cells.Item(495, 9).Value = 'Segunda'
$ws.Cells.Item(495, 10).Value = 830
$ws.Cells.Item(495, 11).Value = 4000
$ws.Cells.Item(495, 12).Value = 5000
$ws.Cells.Item(495, 13).Value = 4723
$ws.Cells.Item(495, 14).Value = '$/docena de matas'
$ws.Cells.Item(495, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(495, 16).Value = 787
$ws.Cells.Item(495, 17).Value = 6
$ws.Cells.Item(495, 18).Value = 'Hortaliza'
